$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1 — copy the format from the neighboring
# header cell (G1, "sum") so it gets the same bold/centered/bordered style,
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# New data value for the "Save" column in row 2
$ws.Range("H2").Value = 0
